$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-08-31 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-01 Friday", 2)

# Update the division problems table (5 data rows x 5 columns, at table rows 1,5,9,13,17)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "82÷4="
$t.Cell(1, 2).Range.Text = "88÷4="
$t.Cell(1, 3).Range.Text = "50÷8="
$t.Cell(1, 4).Range.Text = "79÷8="
$t.Cell(1, 5).Range.Text = "90÷5="

$t.Cell(5, 1).Range.Text = "62÷3="
$t.Cell(5, 2).Range.Text = "54÷5="
$t.Cell(5, 3).Range.Text = "86÷2="
$t.Cell(5, 4).Range.Text = "84÷4="
$t.Cell(5, 5).Range.Text = "48÷3="

$t.Cell(9, 1).Range.Text = "62÷4="
$t.Cell(9, 2).Range.Text = "52÷4="
$t.Cell(9, 3).Range.Text = "34÷2="
$t.Cell(9, 4).Range.Text = "74÷9="
$t.Cell(9, 5).Range.Text = "22÷9="

$t.Cell(13, 1).Range.Text = "53÷3="
$t.Cell(13, 2).Range.Text = "35÷2="
$t.Cell(13, 3).Range.Text = "62÷7="
$t.Cell(13, 4).Range.Text = "97÷4="
$t.Cell(13, 5).Range.Text = "79÷3="

$t.Cell(17, 1).Range.Text = "24÷8="
$t.Cell(17, 2).Range.Text = "66÷7="
$t.Cell(17, 3).Range.Text = "26÷6="
$t.Cell(17, 4).Range.Text = "48÷9="
$t.Cell(17, 5).Range.Text = "38÷8="
